$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) New test-data columns / rows.
#    Cell values are entered in the same order the original author appears to
#    have used (this matters because it controls the order new entries are
#    appended to the shared-strings table).
# ---------------------------------------------------------------------------

# L1 gets a new "id=" marker value and a yellow highlight.
$ws.Range("L1").Value = "id="""
$ws.Range("L1:L2").Interior.Color = 65535

# L4 / C4 / E4 expected-state markers.
$ws.Range("L4").Value = "<NOTEXISTS>"
$ws.Range("C4").Value = "<ENABLED>"
$ws.Range("E4").Value = "<EMPTY>"

# Highlight the whole spacer row (C2:T2) yellow as well.
$ws.Range("C2:T2").Interior.Color = 65535

# Row 1: remaining "id=..." element locator values.
$ws.Range("C1").Value = "id=tricentis_logo"
$ws.Range("D1").Value = "id=visitsupport"
$ws.Range("E1").Value = "id=search_form"
$ws.Range("F1").Value = "id=search_button"
$ws.Range("G1").Value = "id=nav_automobile"
$ws.Range("H1").Value = "id=nav_truck"
$ws.Range("I1").Value = "id=nav_motorcycle"
$ws.Range("J1").Value = "id=nav_camper"
$ws.Range("K1").Value = "id=downloadtrial"
$ws.Range("M1").Value = "id=tricentis_about"
$ws.Range("N1").Value = "id=tricentis_products"
$ws.Range("O1").Value = "id=tricentis_events"
$ws.Range("P1").Value = "id=tricentis_resources"
$ws.Range("Q1").Value = "id=tricentis_services"
$ws.Range("R1").Value = "id=nav_facebook"
$ws.Range("S1").Value = "id=nav_twitter"
$ws.Range("T1").Value = "id=tricentis_services"

# Row 4: remaining <ENABLED> markers (reuse the existing string).
$ws.Range("F4").Value = "<ENABLED>"
$ws.Range("R4").Value = "<ENABLED>"
$ws.Range("S4").Value = "<ENABLED>"
$ws.Range("T4").Value = "<ENABLED>"

# New test rows 6-8: clicking the Truck / Motorcycle / Camper links.
$ws.Range("A6").Value = "Click Truck link"
$ws.Range("B6").Value = "<SET>"
$ws.Range("H6").Value = "X"
$ws.Range("U6").Value = "<NOP>"

$ws.Range("A7").Value = "Click Motorcycle link"
$ws.Range("B7").Value = "<SET>"
$ws.Range("I7").Value = "X"
$ws.Range("U7").Value = "<NOP>"

$ws.Range("A8").Value = "Click Camper link"
$ws.Range("B8").Value = "<SET>"
$ws.Range("J8").Value = "X"
$ws.Range("U8").Value = "<NOP>"

# ---------------------------------------------------------------------------
# 2) Column width adjustments.
#    (XML "width" = COM ColumnWidth + 5/6 in this engine.)
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 27.276041666666668
$ws.Columns.Item(4).ColumnWidth = 13.498697916666666
$ws.Columns.Item(8).ColumnWidth = 12.166666666666666
$ws.Columns.Item(9).ColumnWidth = 17.385416666666668
$ws.Columns.Item(10).ColumnWidth = 14.166666666666666
$ws.Columns.Item(11).ColumnWidth = 15.608072916666666
$ws.Columns.Item(12).ColumnWidth = 11.166666666666666
$ws.Columns.Item(13).ColumnWidth = 16.608072916666668
$ws.Columns.Item(14).ColumnWidth = 18.944010416666668
$ws.Columns.Item(15).ColumnWidth = 19.721354166666668
$ws.Columns.Item(16).ColumnWidth = 18.166666666666668
$ws.Columns.Item(17).ColumnWidth = 16.721354166666668
$ws.Columns.Item(18).ColumnWidth = 14.166666666666666
$ws.Columns.Item(19).ColumnWidth = 12.053385416666666
$ws.Columns.Item(20).ColumnWidth = 16.721354166666668

# ---------------------------------------------------------------------------
# 3) Reposition / resize the logo picture (must happen after the column
#    widths above are set, since the picture's cell anchor is computed from
#    the current column widths / row heights).
# ---------------------------------------------------------------------------
$pic = $ws.Shapes.Item(1)
$pic.Left = 1.8
$pic.Top = 139.2
$pic.Width = 1498.9603937007873

# ---------------------------------------------------------------------------
# 4) Selection state.
# ---------------------------------------------------------------------------
$ws.Range("I7").Select() | Out-Null
